$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January")

# Add the new "type of testing" column (column C) next to the existing names in column B
$ws.Range("C2").Value = "Migration Testing"
$ws.Range("C3").Value = "Automation Testing"
$ws.Range("C4").Value = "API Testing"
$ws.Range("C5").Value = "Performance Testing"

# Size column C to fit its new content (matches the author's best-fit column width)
$ws.Columns.Item(3).ColumnWidth = 15

# Leave the selection on the last filled cell, as in the source workbook
$ws.Range("C5").Select() | Out-Null
